$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift the week's dates forward by 6 days (one week later: 2/6 - 2/12)
$ws.Range("B5").Value = 44598
$ws.Range("C5").Value = 44599
$ws.Range("D5").Value = 44600
$ws.Range("E5").Value = 44601
$ws.Range("F5").Value = 44602
$ws.Range("G5").Value = 44603
$ws.Range("H5").Value = 44604

# Fill in hours logged on Wednesday (column E) for Team Meeting, Sponsor Meeting, TA Meeting
$ws.Range("E8").Value = 1
$ws.Range("E9").Value = 1.5
$ws.Range("E10").Value = 1

# Update the active selection to E12
$ws.Range("E12").Select()
